# Update the "dSF" column (F) values for a handful of rows to reflect
# re-pulled data / recalculated means.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F6").Value = 2
$ws.Range("F8").Value = -3
$ws.Range("F9").Value = -12
$ws.Range("F12").Value = 3
